$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("D", "E", "F", "G", "H", "I", "J")
foreach ($col in $columns) {
    $addr = "${col}2"
    $ws.Range($addr).Value = "unknown"
}

$wb.Save()
